$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for specific rows to match repulled data
$ws.Range("F5").Value = -6
$ws.Range("F7").Value = -3
$ws.Range("F9").Value = -5
$ws.Range("F11").Value = -4
$ws.Range("F14").Value = 3
$ws.Range("F16").Value = 5
$ws.Range("F18").Value = -2
$ws.Range("F19").Value = -10
$ws.Range("F22").Value = 0
